$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final full data set (cluster name, active cases), sorted as in target workbook.
$data = @(
    @("126 Racecourse Road Public Housing Tower Flemington", 10),
    @("139 Highett St Apartment Complex Richmond", 11),
    @("3175 The Bays Aged Care Facility Hastings", 17),
    @("3600 Belvedere Age Care Noble Park Outbreak", 18),
    @("3612 BlueCross Glengowrie Outbreak", 22),
    @("3684 Homestyle Aged Care Langford Grange Cranbourne East Outbreak", 20),
    @("4075 Ferndale Gardens Aged Care Services Bayswater North Outbreak", 16),
    @("Australian Lamb Colac East", 14),
    @("Bread Solutions Braeside Outbreak", 19),
    @("CS Square Caroline Springs Outbreak", 16),
    @("Cedar Meats Australia Brooklyn Outbreak", 11),
    @("Child's Play Early Learning Centre Tarneit", 10),
    @("Embracia Aged Care Reservoir Outbreak", 23),
    @("FedEx Station Melbourne Airport", 10),
    @("Guardian Childcare Caulfield Outbreak", 17),
    @("Inghams Enterprise Somerville Outbreak", 12),
    @("Kool Kidz Childcare Narre Warren", 11),
    @("Lantmannen Unibake Australia Mordialloc", 26),
    @("Nido Early School Ascot Vale", 15),
    @("Nido Early School Glenroy", 15),
    @("Northern Health Northern Hospital Epping Emergency Department Tier 1B", 44),
    @("Northern Health The Northern Hospital Epping", 13),
    @("Oceania Meat Processors Laverton North Outbreak", 16),
    @("Pick It Up Fitness Mulgrave Outbreak", 11),
    @("Robin Hood Inn Drouin West Outbreak", 32),
    @("Social Gathering Warrnambool 28 Sep Outbreak", 17),
    @("St Vincents Hospital Emergency Department Melbourne", 39),
    @("Target Distribution Centre Truganina Outbreak", 21),
    @("The Royal Children's Hospital Melbourne Emergency Department Parkville Tier 1B", 12),
    @("The Toolshed Bar Private Event Noojee", 10),
    @("Turosi Breakwater", 10),
    @("Visy Recycling Springvale", 31),
    @("Werribee Mercy Hospital Emergency Department", 25),
    @("Western Health Sunshine Hospital Emergency Department", 21)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
